# Added viewing of goals
# Highlight the "Create/Read/Update/Delete budget items" user-story
# paragraphs (and the bonus "suggested budget categories" story) green,
# matching how the rest of the user-story list is already highlighted.

$d = $word.ActiveDocument

# Distinctive substrings that uniquely identify each paragraph that must
# be highlighted bright green (wdBrightGreen = 4), both the paragraph
# mark and every run in the paragraph.
$targets = @(
    "create my own budget categories so I can be as specific",
    "enter recurring income which occurs weekly",
    "enter recurring outgo (expense) which occurs weekly",
    "enter single-occurrence income on a particular date",
    "able to enter a single-occurrence expense on a particular date",
    "remove any recurring or single-occurrence income or outgo items",
    "select from suggested budget categories so I can begin budgeting"
)

foreach ($para in $d.Paragraphs) {
    $text = $para.Range.Text
    foreach ($needle in $targets) {
        if ($text -like "*$needle*") {
            $para.Range.Font.HighlightColorIndex = 4
            break
        }
    }
}
